$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the J1 / AVR_ICSP_3x2 row (previously row 4, including the blank
# Id-only rows 2 and 3 above it) and the SW2 / Rotary encoder row (row 12).
$ws.Range("A2:F4").ClearContents()
$ws.Range("A12:F12").ClearContents()

# New manufacturer / manufacturer-part-number data for several components.
$ws.Range("F5").Value = "HRO"
$ws.Range("G5").Value = "K2-1187SQ-D4SW-06"

# Header row: "Supplier and ref" -> "Manufacturer", plus new "Manufacturer
# Part Number" column.
$ws.Range("F1").Value = "Manufacturer"
$ws.Range("G1").Value = "Manufacturer Part Number"

$ws.Range("F7").Value = "Jing Extension of the Electronic Co."
$ws.Range("G7").Value = "920-462A2021D10102"

$ws.Range("F8").Value = "Yangxing Tech"
$ws.Range("G8").Value = "X322516MLB4SI"

$g17 = $ws.Range("G17")
$g17.Value = "MF-MSMF075-2"
$g17Font = $g17.Font
$g17Font.Name = "Inherit"
$g17Font.Size = 10
$g17Font.Color = 0
$g17.VerticalAlignment = -4108
$g17.WrapText = $true

$ws.Range("F17").Value = "Bourns Inc."

$ws.Range("F9").Value = "HRO"
$ws.Range("G9").Value = "TYPE-C-31-M-12"

# New column G needs the same width treatment as column B (23.33203125).
$ws.Columns.Item(7).ColumnWidth = 22.5

# Page setup was touched in the edit session (portrait orientation).
$ws.PageSetup.Orientation = 1

# Reflect the last-used selection from the edit session.
[void]$ws.Range("C12").Select()

Write-Output "done"
